$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC09_Verify_Documented_Savings")

# Insert a new blank row at position 5 (pushes the old rows 5 & 6 down to 6 & 7)
$ws.Rows.Item(5).Insert()

# Re-apply the bordered "data row" look (border on all sides, no fill) that the
# rest of the table uses, so the new row matches its neighbours visually.
$ws.Range("A5:E5").Borders.LineStyle = 1

# New row 5 only carries a single value: the WAIT keyword in column B.
$ws.Range("B5").Value = "WAIT"

# Selection ends up on the freshly entered cell.
$ws.Range("B5").Select()
